# Applies the "Updated symbol list" GitHub Actions data refresh:
# refreshed Price (col D) / Volume 1h (col E) figures for most coins,
# and swapped the BOLO / CoinbaseStockToken rows (47 & 48) with updated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.52"
$ws.Range("E2").Value = "'1.75%"
$ws.Range("D3").Value = "'31.14"
$ws.Range("E3").Value = "'-0.33%"
$ws.Range("D4").Value = "'5.133"
$ws.Range("E4").Value = "'0.85%"
$ws.Range("D5").Value = "'0.08120"
$ws.Range("E5").Value = "'10.27%"
$ws.Range("D6").Value = "'2.507"
$ws.Range("E6").Value = "'51.73%"
$ws.Range("D7").Value = "'7.843"
$ws.Range("E7").Value = "'2.05%"
$ws.Range("D8").Value = "'3.843"
$ws.Range("E8").Value = "'2.17%"
$ws.Range("D9").Value = "'0.9116"
$ws.Range("E9").Value = "'-1.32%"
$ws.Range("E10").Value = "'2.31%"
$ws.Range("D11").Value = "'0.07258"
$ws.Range("E11").Value = "'1.90%"
$ws.Range("D12").Value = "'0.08019"
$ws.Range("E12").Value = "'1.98%"
$ws.Range("D13").Value = "'0.03026"
$ws.Range("E13").Value = "'0.82%"
$ws.Range("D14").Value = "'0.09971"
$ws.Range("E14").Value = "'0.79%"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("E15").Value = "'0.32%"
$ws.Range("D16").Value = "'0.005982"
$ws.Range("E16").Value = "'-4.40%"
$ws.Range("E17").Value = "'1.13%"
$ws.Range("D18").Value = "'2.241"
$ws.Range("E18").Value = "'0.55%"
$ws.Range("D19").Value = "'0.3309"
$ws.Range("E19").Value = "'0.93%"
$ws.Range("D20").Value = "'0.1348"
$ws.Range("E20").Value = "'-0.08%"
$ws.Range("D21").Value = "'4.595"
$ws.Range("E21").Value = "'0.70%"
$ws.Range("D22").Value = "'0.1604"
$ws.Range("E22").Value = "'3.27%"
$ws.Range("D23").Value = "'0.04591"
$ws.Range("D24").Value = "'0.001260"
$ws.Range("E24").Value = "'3.34%"
$ws.Range("D25").Value = "'0.004445"
$ws.Range("E25").Value = "'0.61%"
$ws.Range("D26").Value = "'0.0001182"
$ws.Range("E26").Value = "'-9.17%"
$ws.Range("D27").Value = "'0.0003438"
$ws.Range("E27").Value = "'83.08%"
$ws.Range("D39").Value = "'0.01812"
$ws.Range("E39").Value = "'9.37%"
$ws.Range("D40").Value = "'0.04537"
$ws.Range("E40").Value = "'3.22%"
$ws.Range("D41").Value = "'0.007224"
$ws.Range("E41").Value = "'2.02%"
$ws.Range("D42").Value = "'0.1343"
$ws.Range("E42").Value = "'1.36%"
$ws.Range("D43").Value = "'0.002176"
$ws.Range("E43").Value = "'3.53%"
$ws.Range("D44").Value = "'0.01071"
$ws.Range("E44").Value = "'-2.73%"
$ws.Range("D45").Value = "'0.00006262"
$ws.Range("E45").Value = "'4.10%"
$ws.Range("E46").Value = "'0.24%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.8206"
$ws.Range("E47").Value = "'15.31%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.006658"
$ws.Range("E48").Value = "'-39.55%"
$ws.Range("D49").Value = "'0.00002105"
$ws.Range("E49").Value = "'0.24%"
$ws.Range("D50").Value = "'0.0002005"
$ws.Range("E50").Value = "'0.31%"
